# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
# Match the header style used by the existing header cells (e.g. AC1) first,
# then set the text for each new header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every data row (2 through 58)
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # column AD
    $ws.Cells.Item($r, 31).Value = 84   # column AE
    $ws.Cells.Item($r, 32).Value = 0    # column AF
}
